$d = $word.ActiveDocument
$d.Content.Find.Execute("rotated incrementally 360° around fruit", $true, $false, $false, $false, $false, $true, 1, $false, "rotated incrementally 360° around the fruit", 2)
